$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '57.666.28'
$ws.Range("E2").Value = '  +2.61%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.321.64'
$ws.Range("E3").Value = '  +0.80%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '539.56'
$ws.Range("E5").Value = '  +4.40%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.53'
$ws.Range("E6").Value = '  +2.83%  '

$ws.Range("E7").Value = '  +0.30%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.559'
$ws.Range("E8").Value = '  +5.54%  '

$ws.Range("E9").Value = '  +1.45%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.45'
$ws.Range("E10").Value = '  +3.89%  '

$ws.Range("E11").Value = '  -0.26%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.357'
$ws.Range("E12").Value = '  +6.11%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '23.62'
$ws.Range("E13").Value = '  +1.77%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.734.99'
$ws.Range("E14").Value = '  +0.67%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '57.630.67'
$ws.Range("E15").Value = '  +2.64%  '

$ws.Range("E16").Value = '  +0.99%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.329.00'
$ws.Range("E17").Value = '  +0.96%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.58'
$ws.Range("E18").Value = '  +2.63%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '331.78'
$ws.Range("E19").Value = '  +0.67%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.23'
$ws.Range("E20").Value = '  +2.64%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.64'
$ws.Range("E21").Value = '  -0.53%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.997'
$ws.Range("E22").Value = '  -0.18%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.58'
$ws.Range("E23").Value = '  +0.76%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '62.04'
$ws.Range("E24").Value = '  +1.92%  '

$ws.Range("E25").Value = '  +2.74%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  +0.48%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.43'
$ws.Range("E27").Value = '  -1.70%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.40'
$ws.Range("E28").Value = '  +6.26%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.75'
$ws.Range("E29").Value = '  +4.00%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '170.51'
$ws.Range("E30").Value = '  +1.57%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0₃0727'
$ws.Range("E31").Value = '  +2.36%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.09'
$ws.Range("E32").Value = '  +0.20%  '

$ws.Range("E33").Value = '  +16.66%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '18.35'
$ws.Range("E34").Value = '  +0.86%  '

$ws.Range("E35").Value = '  +0.03%  '

$ws.Range("E36").Value = '  +0.36%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.15'
$ws.Range("E37").Value = '  +7.05%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.24'
$ws.Range("E38").Value = '  +0.88%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.61'
$ws.Range("E39").Value = '  +3.50%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '39.07'
$ws.Range("E40").Value = '  +1.38%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '143.93'
$ws.Range("E41").Value = '  -2.62%  '

$ws.Range("E42").Value = '  +0.80%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.61'
$ws.Range("E43").Value = '  +1.80%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '283.53'
$ws.Range("E44").Value = '  -0.05%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0934'
$ws.Range("E45").Value = '  +1.03%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '19.01'
$ws.Range("E46").Value = '  +4.90%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0500'
$ws.Range("E47").Value = '  +0.99%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.557'
$ws.Range("E48").Value = '  +0.53%  '

$ws.Range("E49").Value = '  +2.33%  '

$ws.Range("E50").Value = '  +1.18%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '17.40'
$ws.Range("E51").Value = '  +1.81%  '
